$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shared across the Overview sheet (zh-cn/de-de status
#    columns) and the per-locale sheets' "Status" column, for both rows.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: fill in "Latest Target File" handback info for both rows
#    (file a.md and b.md both report back the same handed-back xlf), and
#    update the handback datetime.
# ---------------------------------------------------------------------------
$zhXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhDateTime = "2016-08-31 16:41:53"

$wsZh.Range("J2").Value = $zhXlf
$wsZh.Range("K2").Value = $zhDateTime
$wsZh.Range("J3").Value = $zhXlf
$wsZh.Range("K3").Value = $zhDateTime

# ---------------------------------------------------------------------------
# 3) de-de sheet: same as above, with the de-de xlf + its own datetime.
# ---------------------------------------------------------------------------
$deXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deDateTime = "2016-08-31 16:42:11"

$wsDe.Range("J2").Value = $deXlf
$wsDe.Range("K2").Value = $deDateTime
$wsDe.Range("J3").Value = $deXlf
$wsDe.Range("K3").Value = $deDateTime

# ---------------------------------------------------------------------------
# 4) Add "a.md" hyperlinks in the "Latest Handback File" column (I) for both
#    rows on both locale sheets. Rebuild each sheet's whole hyperlink
#    collection (existing + new) so relationship ids come out in the same
#    left-to-right, top-to-bottom order as the handed-back workbook.
# ---------------------------------------------------------------------------
$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c37f94a981dbd9add50f1420a525f814d04f7d50/e2e/a.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c37f94a981dbd9add50f1420a525f814d04f7d50/e2e/b.md"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlA, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlB, "", "", "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlA, "", "", "a.md")

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlA, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlB, "", "", "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlA, "", "", "a.md")

# ---------------------------------------------------------------------------
# 5) Column width adjustments caused by the wider "Status" text / new
#    hyperlink file name text (what Excel's own AutoFit would produce for
#    these columns after the content above was written).
#    ColumnWidth is expressed in characters; Excel quantizes the stored
#    width to whole pixels (characters*6 rounded, +5 padding, /6), so the
#    inputs below are chosen to land on the nearest achievable width to the
#    handed-back file's 29.9777.../40 column widths.
# ---------------------------------------------------------------------------
$wsOverview.Columns("E").ColumnWidth = 29.166666666666668
$wsOverview.Columns("F").ColumnWidth = 29.166666666666668

$wsZh.Columns("C").ColumnWidth = 29.166666666666668
$wsZh.Columns("J").ColumnWidth = 39.166666666666664

$wsDe.Columns("C").ColumnWidth = 29.166666666666668
$wsDe.Columns("J").ColumnWidth = 39.166666666666664
